# B6-PowerPoint.pptx edit: Wed, May 20, 2020  1:05:09 AM
#
# The commit swaps the table style applied to the three tables in the
# deck (slides 14, 15 and 16) from the custom "Table_0" style
# ({03D36A11-2327-4319-8E84-E945BA12BCCE}) to PowerPoint's built-in
# "No Style, No Grid" table style ({DE67CE6F-F196-483E-B238-6C4AA5BD93C1}).
#
# (The commit also records embedding the Tahoma font used elsewhere in
# the deck via File > Options > Save > "Embed fonts in the file" -- that
# is a Save-dialog-only feature with no COM/VBA automation surface in
# PowerPoint, so it cannot be reproduced through the object model.)

$p = $ppt.ActivePresentation

$newStyleId = "{DE67CE6F-F196-483E-B238-6C4AA5BD93C1}"

$tableSlideIndexes = 14, 15, 16

foreach ($slideIndex in $tableSlideIndexes) {
    $slide = $p.Slides.Item($slideIndex)
    for ($shapeIndex = 1; $shapeIndex -le $slide.Shapes.Count; $shapeIndex++) {
        $shape = $slide.Shapes.Item($shapeIndex)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId, $true)
        }
    }
}

Write-Output "table styles updated"
